# Updates the cryptos worksheet with refreshed price/volume data
# (commit: "Updated cryptos list ... with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: rows 41/42 swap places (Kaspa <-> Bittensor) and are updated together
# with every other changed B/C/E cell in the $otherUpdates table below.

# --- Price ("D" column) updates -------------------------------------------
# These values must stay plain text (matching the source XML's inlineStr
# cells), so force a Text number format before assigning, then restore the
# default "Normal" style so no stray formatting differences are introduced.
$priceUpdates = @(
    @{Cell='D2'; Value='59.877.33'}
    @{Cell='D3'; Value='2.964.16'}
    @{Cell='D5'; Value='574.35'}
    @{Cell='D6'; Value='124.65'}
    @{Cell='D8'; Value='2.959.72'}
    @{Cell='D9'; Value='0.501'}
    @{Cell='D12'; Value='0.438'}
    @{Cell='D14'; Value='32.38'}
    @{Cell='D16'; Value='3.465.83'}
    @{Cell='D17'; Value='2.967.29'}
    @{Cell='D18'; Value='59.914.06'}
    @{Cell='D19'; Value='6.18'}
    @{Cell='D20'; Value='431.14'}
    @{Cell='D21'; Value='13.05'}
    @{Cell='D22'; Value='0.659'}
    @{Cell='D23'; Value='6.98'}
    @{Cell='D24'; Value='12.64'}
    @{Cell='D25'; Value='78.89'}
    @{Cell='D27'; Value='1.00'}
    @{Cell='D28'; Value='2.52'}
    @{Cell='D29'; Value='7.21'}
    @{Cell='D30'; Value='1.87'}
    @{Cell='D31'; Value='6.10'}
    @{Cell='D32'; Value='25.21'}
    @{Cell='D33'; Value='0.0929'}
    @{Cell='D34'; Value='2.17'}
    @{Cell='D35'; Value='0.944'}
    @{Cell='D36'; Value='5.57'}
    @{Cell='D37'; Value='49.35'}
    @{Cell='D39'; Value='7.92'}
    @{Cell='D41'; Value='381.63'}
    @{Cell='D42'; Value='0.108'}
    @{Cell='D44'; Value='2.621.07'}
    @{Cell='D47'; Value='118.58'}
    @{Cell='D50'; Value='23.25'}
    @{Cell='D51'; Value='31.17'}
)

foreach ($item in $priceUpdates) {
    $cell = $ws.Range($item.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# --- Volume(1h) ("E" column, and remaining B/C) updates --------------------
$otherUpdates = @(
    @{Cell='E2'; Value='  -5.03%  '}
    @{Cell='E3'; Value='  -6.65%  '}
    @{Cell='E4'; Value='  +0.08%  '}
    @{Cell='E5'; Value='  -2.92%  '}
    @{Cell='E6'; Value='  -7.61%  '}
    @{Cell='E7'; Value='  +0.09%  '}
    @{Cell='E8'; Value='  -6.82%  '}
    @{Cell='E9'; Value='  -2.87%  '}
    @{Cell='E10'; Value='  -6.40%  '}
    @{Cell='E11'; Value='  -2.91%  '}
    @{Cell='E12'; Value='  -3.31%  '}
    @{Cell='E13'; Value='  -6.71%  '}
    @{Cell='E14'; Value='  -6.15%  '}
    @{Cell='E15'; Value='  -0.38%  '}
    @{Cell='E16'; Value='  -6.35%  '}
    @{Cell='E17'; Value='  -6.59%  '}
    @{Cell='E18'; Value='  -4.88%  '}
    @{Cell='E19'; Value='  -5.62%  '}
    @{Cell='E20'; Value='  -6.38%  '}
    @{Cell='E21'; Value='  -6.85%  '}
    @{Cell='E22'; Value='  -5.53%  '}
    @{Cell='E23'; Value='  -8.24%  '}
    @{Cell='E24'; Value='  -4.81%  '}
    @{Cell='E25'; Value='  -4.26%  '}
    @{Cell='E26'; Value='  +0.16%  '}
    @{Cell='E27'; Value='  -0.03%  '}
    @{Cell='E28'; Value='  -5.58%  '}
    @{Cell='E29'; Value='  -5.64%  '}
    @{Cell='E30'; Value='  -7.81%  '}
    @{Cell='E31'; Value='  -9.31%  '}
    @{Cell='E32'; Value='  -7.27%  '}
    @{Cell='E33'; Value='  -8.78%  '}
    @{Cell='E34'; Value='  -8.30%  '}
    @{Cell='E35'; Value='  -8.06%  '}
    @{Cell='E36'; Value='  -3.93%  '}
    @{Cell='E37'; Value='  -3.89%  '}
    @{Cell='E38'; Value='  -7.78%  '}
    @{Cell='E39'; Value='  -2.02%  '}
    @{Cell='E40'; Value='  -7.68%  '}
    @{Cell='B41'; Value='Bittensor'}
    @{Cell='C41'; Value='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'}
    @{Cell='E41'; Value='  -5.56%  '}
    @{Cell='B42'; Value='Kaspa'}
    @{Cell='C42'; Value='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'}
    @{Cell='E42'; Value='  -2.92%  '}
    @{Cell='E43'; Value='  -6.95%  '}
    @{Cell='E44'; Value='  -6.89%  '}
    @{Cell='E45'; Value='  +0.06%  '}
    @{Cell='E46'; Value='  -6.56%  '}
    @{Cell='E47'; Value='  -4.68%  '}
    @{Cell='E48'; Value='  -6.53%  '}
    @{Cell='E49'; Value='  -4.33%  '}
    @{Cell='E50'; Value='  -7.31%  '}
    @{Cell='E51'; Value='  -10.42%  '}
)

foreach ($item in $otherUpdates) {
    $ws.Range($item.Cell).Value = $item.Value
}
